$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 19854.23
$ws.Range("J17").Value = 19854.23
$ws.Range("L17").Value = 59562.69
$ws.Range("N17").Value = -59898.69
$ws.Range("H38").Value = 728.6875
$ws.Range("I38").Value = 78.09091
$ws.Range("J38").Value = 2160.0
$ws.Range("K38").Value = 234.27273
$ws.Range("L38").Value = 6480.0
$ws.Range("M38").Value = 137.72727
$ws.Range("N38").Value = -7224.0
$ws.Range("H40").Value = 1507.95
$ws.Range("I40").Value = 1373.7693
$ws.Range("J40").Value = 1757.1428
$ws.Range("K40").Value = 1373.7693
$ws.Range("L40").Value = 1757.1428
$ws.Range("M40").Value = -1198.7693
$ws.Range("N40").Value = -2107.1428
$ws.Range("H46").Value = 500001250.0
$ws.Range("I46").Value = 0.0
$ws.Range("J46").Value = 500001250.0
$ws.Range("K46").Value = 0.0
$ws.Range("L46").Value = 1500003750.0
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -1500003988.0
$ws.Range("H60").Value = 500001250.0
$ws.Range("I60").Value = 0.0
$ws.Range("J60").Value = 500001250.0
$ws.Range("K60").Value = 0.0
$ws.Range("L60").Value = 1500003750.0
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = -1500004718.0
$ws.Range("H76").Value = 3975.3845
$ws.Range("I76").Value = 3886.6667
$ws.Range("J76").Value = 4175.0
$ws.Range("K76").Value = 3886.6667
$ws.Range("L76").Value = 4175.0
$ws.Range("M76").Value = -3571.6667
$ws.Range("N76").Value = -4805.0
$ws.Range("H79").Value = 3975.3845
$ws.Range("I79").Value = 3886.6667
$ws.Range("J79").Value = 4175.0
$ws.Range("K79").Value = 3886.6667
$ws.Range("L79").Value = 4175.0
$ws.Range("M79").Value = -2794.6667
$ws.Range("N79").Value = -6359.0
$ws.Range("H86").Value = 1268.2222
$ws.Range("I86").Value = 1361.5
$ws.Range("K86").Value = 1361.5
$ws.Range("M86").Value = -238.5
$ws.Range("H89").Value = 1268.2222
$ws.Range("I89").Value = 1361.5
$ws.Range("K89").Value = 6807.5
$ws.Range("M89").Value = -1191.5
$ws.Range("H113").Value = 3300.0
$ws.Range("I113").Value = 2000.0
$ws.Range("J113").Value = 3625.0
$ws.Range("K113").Value = 2000.0
$ws.Range("L113").Value = 3625.0
$ws.Range("M113").Value = 1254.0
$ws.Range("N113").Value = -10133.0
$ws.Range("H129").Value = 1010.1579
$ws.Range("J129").Value = 1075.0312
$ws.Range("L129").Value = 3225.0936
$ws.Range("N129").Value = -13225.0936
$ws.Range("H137").Value = 1718.9565
$ws.Range("I137").Value = 1563.0834
$ws.Range("J137").Value = 1889.0
$ws.Range("K137").Value = 4689.2502
$ws.Range("L137").Value = 5667.0
$ws.Range("M137").Value = -2139.2502
$ws.Range("N137").Value = -10767.0
$ws.Range("H138").Value = 4317.3877
$ws.Range("I138").Value = 1380.8572
$ws.Range("J138").Value = 8232.762
$ws.Range("K138").Value = 4142.571599999999
$ws.Range("L138").Value = 24698.286
$ws.Range("M138").Value = 997.4284000000007
$ws.Range("N138").Value = -34978.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11592.107
$ws.Range("I32").Value = 14268.651
$ws.Range("J32").Value = 2738.923
$ws.Range("K32").Value = 14268.651
$ws.Range("L32").Value = 2738.923
$ws.Range("M32").Value = -13981.651
$ws.Range("N32").Value = -3312.923
$ws.Range("H61").Value = 977.4
$ws.Range("I61").Value = 843.14636
$ws.Range("J61").Value = 1370.5714
$ws.Range("K61").Value = 843.14636
$ws.Range("L61").Value = 1370.5714
$ws.Range("M61").Value = -631.14636
$ws.Range("N61").Value = -1794.5714
$ws.Range("H74").Value = 609.71875
$ws.Range("I74").Value = 592.7451
$ws.Range("J74").Value = 676.3077
$ws.Range("K74").Value = 592.7451
$ws.Range("L74").Value = 676.3077
$ws.Range("M74").Value = 281.2549
$ws.Range("N74").Value = -2424.3077
$ws.Range("H77").Value = 609.71875
$ws.Range("I77").Value = 592.7451
$ws.Range("J77").Value = 676.3077
$ws.Range("K77").Value = 2963.7255
$ws.Range("L77").Value = 3381.5385
$ws.Range("M77").Value = 1404.2745
$ws.Range("N77").Value = -12117.5385
$ws.Range("H110").Value = 893.2222
$ws.Range("I110").Value = 907.65515
$ws.Range("J110").Value = 833.4286
$ws.Range("K110").Value = 907.65515
$ws.Range("L110").Value = 833.4286
$ws.Range("M110").Value = 1137.34485
$ws.Range("N110").Value = -4923.4286
$ws.Range("H122").Value = 1496.5883
$ws.Range("I122").Value = 1296.1333
$ws.Range("J122").Value = 3000.0
$ws.Range("K122").Value = 3888.3999
$ws.Range("L122").Value = 9000.0
$ws.Range("M122").Value = -1438.3999
$ws.Range("N122").Value = -13900.0
$ws.Range("H132").Value = 2131.279
$ws.Range("I132").Value = 1513.2354
$ws.Range("J132").Value = 2535.3845
$ws.Range("K132").Value = 4539.706200000001
$ws.Range("L132").Value = 7606.1535
$ws.Range("M132").Value = -2009.706200000001
$ws.Range("N132").Value = -12666.1535
$ws.Range("H133").Value = 84468.37
$ws.Range("J133").Value = 84468.37
$ws.Range("L133").Value = 84468.37
$ws.Range("N133").Value = -89528.37
$ws.Range("H136").Value = 977.4
$ws.Range("I136").Value = 843.14636
$ws.Range("J136").Value = 1370.5714
$ws.Range("K136").Value = 2529.43908
$ws.Range("L136").Value = 4111.7142
$ws.Range("M136").Value = 20.5609199999999
$ws.Range("N136").Value = -9211.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3342.3125
$ws.Range("I105").Value = 3300.6
$ws.Range("J105").Value = 3411.8333
$ws.Range("K105").Value = 3300.6
$ws.Range("L105").Value = 3411.8333
$ws.Range("M105").Value = -1553.6
$ws.Range("N105").Value = -6905.8333
$ws.Range("H132").Value = 52841.215
$ws.Range("J132").Value = 52841.215
$ws.Range("L132").Value = 52841.215
$ws.Range("N132").Value = -62961.215
$ws.Range("H133").Value = 59444.832
$ws.Range("J133").Value = 65592.0
$ws.Range("L133").Value = 65592.0
$ws.Range("N133").Value = -75712.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2433.2903
$ws.Range("I31").Value = 1596.762
$ws.Range("J31").Value = 4190.0
$ws.Range("K31").Value = 1596.762
$ws.Range("L31").Value = 4190.0
$ws.Range("M31").Value = -1301.762
$ws.Range("N31").Value = -4780.0
$ws.Range("H34").Value = 2433.2903
$ws.Range("I34").Value = 1596.762
$ws.Range("J34").Value = 4190.0
$ws.Range("K34").Value = 1596.762
$ws.Range("L34").Value = 4190.0
$ws.Range("M34").Value = -1394.762
$ws.Range("N34").Value = -4594.0
$ws.Range("H99").Value = 3553.0
$ws.Range("I99").Value = 4102.4
$ws.Range("J99").Value = 2768.1428
$ws.Range("K99").Value = 4102.4
$ws.Range("L99").Value = 2768.1428
$ws.Range("M99").Value = -2604.4
$ws.Range("N99").Value = -5764.1428
$ws.Range("H126").Value = 3553.0
$ws.Range("I126").Value = 4102.4
$ws.Range("J126").Value = 2768.1428
$ws.Range("K126").Value = 12307.2
$ws.Range("L126").Value = 8304.4284
$ws.Range("M126").Value = -9837.199999999999
$ws.Range("N126").Value = -13244.4284
$ws.Range("H134").Value = 1052.8674
$ws.Range("I134").Value = 946.08954
$ws.Range("J134").Value = 1500.0
$ws.Range("K134").Value = 2838.26862
$ws.Range("L134").Value = 4500.0
$ws.Range("M134").Value = -303.2686200000003
$ws.Range("N134").Value = -9570.0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8480.1
$ws.Range("I56").Value = 8480.1
$ws.Range("K56").Value = 8480.1
$ws.Range("M56").Value = -7950.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 36000.0
$ws.Range("I82").Value = 0.0
$ws.Range("K82").Value = 0.0
$ws.Range("M82").Value = $null
$ws.Range("H85").Value = 36000.0
$ws.Range("I85").Value = 0.0
$ws.Range("K85").Value = 0.0
$ws.Range("M85").Value = $null
$ws.Range("H122").Value = 2911.5557
$ws.Range("I122").Value = 2149.3333
$ws.Range("J122").Value = 4436.0
$ws.Range("K122").Value = 6447.999899999999
$ws.Range("L122").Value = 13308.0
$ws.Range("M122").Value = -3997.999899999999
$ws.Range("N122").Value = -18208.0
$ws.Range("H126").Value = 1611.7
$ws.Range("I126").Value = 1235.2222
$ws.Range("J126").Value = 5000.0
$ws.Range("K126").Value = 3705.6666
$ws.Range("L126").Value = 15000.0
$ws.Range("M126").Value = -1235.6666
$ws.Range("N126").Value = -19940.0
$ws.Range("H132").Value = 2121.1428
$ws.Range("I132").Value = 1220.2222
$ws.Range("J132").Value = 3742.8
$ws.Range("K132").Value = 3660.6666
$ws.Range("L132").Value = 11228.4
$ws.Range("M132").Value = -1130.6666
$ws.Range("N132").Value = -16288.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4777.067
$ws.Range("I7").Value = 3233.4443
$ws.Range("J7").Value = 7092.5
$ws.Range("K7").Value = 3233.4443
$ws.Range("L7").Value = 7092.5
$ws.Range("M7").Value = -3121.4443
$ws.Range("N7").Value = -7316.5
$ws.Range("H126").Value = 4777.067
$ws.Range("I126").Value = 3233.4443
$ws.Range("J126").Value = 7092.5
$ws.Range("K126").Value = 9700.332900000001
$ws.Range("L126").Value = 21277.5
$ws.Range("M126").Value = -7230.332900000001
$ws.Range("N126").Value = -26217.5
$ws.Range("H136").Value = 2876.2144
$ws.Range("I136").Value = 3100.1155
$ws.Range("J136").Value = 2229.389
$ws.Range("K136").Value = 9300.3465
$ws.Range("L136").Value = 6688.167
$ws.Range("M136").Value = -6750.3465
$ws.Range("N136").Value = -11788.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1588.5918
$ws.Range("I136").Value = 1228.8605
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 3686.5815
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -1136.5815
$ws.Range("N136").Value = -17599.9995
